$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sval data (regenerated to filter save games), rows 2-8, columns B-G
$data = @{
    2 = @{ B = 3.286832544864788;    C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; F = 0; G = 6.189590430959694 }
    3 = @{ B = 3.286832544864788;    C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; F = 1; G = 6.189590430959694 }
    4 = @{ B = 0.01293466051926884; C = 1.655778082260271;    D = 0.1494219747398047; E = 0.4942365360607697; F = 1; G = 2.312371253580114 }
    5 = @{ B = 3.286832544864788;    C = 0.04071648406533734;  D = 0.7527432677738641; E = 0.4942365360607697; F = 0; G = 4.574528832764759 }
    6 = @{ B = 0.0006408296065709695; C = 0.00006240767534437808; D = 0.1494219747398047; E = 0.4942365360607697; F = 0; G = 0.6443617480824897 }
    7 = @{ B = 0.1190320826869504;   C = 0.306821227259698;    D = 0.7527432677738641; E = 10.19245300693656;  F = 0; G = 11.37104958465707 }
    8 = @{ B = 0.6606524410359556;   C = 1.655778082260271;    D = 0.1494219747398047; E = 0.4942365360607697; F = 1; G = 2.960089034096801 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
